$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to Text format
# first, otherwise Excel auto-converts them to numeric cells (losing the exact
# literal string representation, e.g. trailing zeros like "62.30" -> 62.3).
$textCells = @("D5", "D7", "D12", "D14", "D15", "D16", "D20", "D22", "D24", "D25", "D26", "D27", "D29", "D31", "D33", "D35", "D36", "D40", "D41", "D46", "D47", "D49")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "38.761.92"
$ws.Range("E2").Value = "  +0.82%  "
$ws.Range("D3").Value = "2.102.96"
$ws.Range("E3").Value = "  +0.64%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "227.41"
$ws.Range("E5").Value = "  -0.41%  "
$ws.Range("E6").Value = "  +0.67%  "
$ws.Range("D7").Value = "62.30"
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  +2.33%  "
$ws.Range("E10").Value = "  +0.41%  "
$ws.Range("E11").Value = "  -1.24%  "
$ws.Range("D12").Value = "15.72"
$ws.Range("E12").Value = "  +6.03%  "
$ws.Range("D13").Value = "2.416.82"
$ws.Range("E13").Value = "  +0.63%  "
$ws.Range("D14").Value = "22.01"
$ws.Range("E14").Value = "  -1.51%  "
$ws.Range("D15").Value = "0.812"
$ws.Range("E15").Value = "  +3.57%  "
$ws.Range("D16").Value = "5.54"
$ws.Range("E16").Value = "  +1.70%  "
$ws.Range("D17").Value = "2.096.95"
$ws.Range("E17").Value = "  +0.27%  "
$ws.Range("D18").Value = "38.786.68"
$ws.Range("E18").Value = "  +1.06%  "
$ws.Range("E19").Value = "  +1.15%  "
$ws.Range("D20").Value = "71.58"
$ws.Range("E20").Value = "  +0.89%  "
$ws.Range("D21").Value = "0.0₃0842"
$ws.Range("E21").Value = "  +1.08%  "
$ws.Range("D22").Value = "228.44"
$ws.Range("E22").Value = "  +1.31%  "
$ws.Range("D24").Value = "2.36"
$ws.Range("E24").Value = "  -3.51%  "
$ws.Range("D25").Value = "2.32"
$ws.Range("E25").Value = "  +0.30%  "
$ws.Range("D26").Value = "9.68"
$ws.Range("E26").Value = "  +2.62%  "
$ws.Range("D27").Value = "172.34"
$ws.Range("E27").Value = "  +1.77%  "
$ws.Range("E28").Value = "  +1.78%  "
$ws.Range("D29").Value = "1.42"
$ws.Range("E29").Value = "  +4.82%  "
$ws.Range("E30").Value = "  +1.47%  "
$ws.Range("D31").Value = "2.54"
$ws.Range("E31").Value = "  +7.97%  "
$ws.Range("E32").Value = "  +0.40%  "
$ws.Range("D33").Value = "4.55"
$ws.Range("E33").Value = "  +0.51%  "
$ws.Range("E34").Value = "  -0.53%  "
$ws.Range("D35").Value = "7.08"
$ws.Range("E35").Value = "  +10.49%  "
$ws.Range("D36").Value = "0.0617"
$ws.Range("E36").Value = "  +1.76%  "
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("E38").Value = "  -0.70%  "
$ws.Range("E39").Value = "  +0.21%  "
$ws.Range("D40").Value = "18.08"
$ws.Range("E40").Value = "  -1.32%  "
$ws.Range("D41").Value = "102.64"
$ws.Range("E41").Value = "  +2.82%  "
$ws.Range("D43").Value = "1.527.38"
$ws.Range("E43").Value = "  -0.75%  "
$ws.Range("E44").Value = "  +8.30%  "
$ws.Range("E45").Value = "  -0.87%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "7.81"
$ws.Range("E46").Value = "  -0.11%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "0.0918"
$ws.Range("E47").Value = "  -2.25%  "
$ws.Range("E48").Value = "  +4.33%  "
$ws.Range("D49").Value = "4.13"
$ws.Range("E49").Value = "  -0.91%  "
$ws.Range("E50").Value = "  -0.55%  "
$ws.Range("D51").Value = "2.302.96"
$ws.Range("E51").Value = "  +0.65%  "
